# diary update wenchia 0303
# Fill in three new diary entries (rows 50-52) that replace the blank
# placeholder rows, matching the style used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DiaryRow($Row, $TemplateRow, $Date, $Time, $Participants, $Goal, $Achievements, $Reflection, $Mood) {
    # Copy the cell formatting (styles) from the template row so the new
    # row reuses the existing style indices instead of creating new ones.
    $cols = @("A", "B", "C", "D", "E", "F", "G")
    foreach ($col in $cols) {
        $ws.Range("$col$TemplateRow").Copy()
        $ws.Range("$col$Row").PasteSpecial(-4122)
    }

    $ws.Cells.Item($Row, 1).Value = $Date
    $ws.Cells.Item($Row, 2).Value = $Time
    $ws.Cells.Item($Row, 3).Value = $Participants
    $ws.Cells.Item($Row, 4).Value = $Goal
    $ws.Cells.Item($Row, 5).Value = $Achievements
    $ws.Cells.Item($Row, 6).Value = $Reflection
    $ws.Cells.Item($Row, 7).Value = $Mood
}

# Row 50 - 2020-02-29
Set-DiaryRow 50 49 43890 "11:00 - 12:00" "Me" "To learn observer pattern and pracrice it with an actul code" "Sucessfully made two small examples and implented with observer pattern" "The example code from the youtube website is not that practical, so I implemented with another example by using event listener. This is more like we will code in the real world situatin." "Feel good!"

# Row 51 - 2020-03-01
Set-DiaryRow 51 49 43891 "13:00 - 14:00" "Me" "To learn decorator pattern and pracrice it with an actul code" "Sucessfully made two small examples and implented with decorator pattern" "Again, the example code from te Youtube isn't good as other resource from other website. Thus, I implemented two examples. The other one is clearer than the previous example." "Feel good!"

# Row 52 - 2020-03-02 (note: the Reflection column (F52) keeps its
# pre-existing placeholder style (italic 11pt, no border) rather than
# being reformatted like D52/E52 - matches the source row exactly)
$cols52 = @("A", "B", "C", "D", "E", "G")
foreach ($col in $cols52) {
    $ws.Range("${col}49").Copy()
    $ws.Range("${col}52").PasteSpecial(-4122)
}

$ws.Cells.Item(52, 1).Value = 43892
$ws.Cells.Item(52, 2).Value = "17:00 - 18:30"
$ws.Cells.Item(52, 3).Value = "Me"
$ws.Cells.Item(52, 4).Value = "To learn factory and abstract facory pattern and pracrice them with an actul code"
$ws.Cells.Item(52, 5).Value = "Sucessfully made two small examples and implented with factory and abstract factory patterns"
$ws.Cells.Item(52, 7).Value = "Feel useful!"
$ws.Cells.Item(52, 6).Value = "This time, I implemented with the pattern by myself and made up some real world situations. Factory patterns are useful and thet often implemented by many application. I'm not unfamiliar with them."

# Row heights to match the wrapped-text content (row 50 has a single
# shorter reflection line; rows 51/52 have longer ones).
$ws.Rows.Item(50).RowHeight = 85
$ws.Rows.Item(51).RowHeight = 102
$ws.Rows.Item(52).RowHeight = 102

# Update the sheet view / selection to match where the author left off
# editing.
$ws.Application.ActiveWindow.ScrollRow = 48
$ws.Range("F53").Select()
